$d = $word.ActiveDocument

$replacements = @(
    @{old = "32÷8="; new = "93÷3="},
    @{old = "55÷4="; new = "65÷7="},
    @{old = "16÷7="; new = "15÷2="},
    @{old = "17÷5="; new = "49÷6="},
    @{old = "44÷9="; new = "93÷6="},
    @{old = "73÷4="; new = "49÷4="},
    @{old = "40÷3="; new = "65÷7="},
    @{old = "37÷2="; new = "92÷7="},
    @{old = "44÷3="; new = "79÷2="},
    @{old = "22÷4="; new = "27÷9="},
    @{old = "56÷4="; new = "65÷6="},
    @{old = "77÷5="; new = "46÷6="},
    @{old = "57÷6="; new = "43÷9="},
    @{old = "97÷9="; new = "35÷3="},
    @{old = "59÷4="; new = "30÷6="},
    @{old = "90÷5="; new = "29÷8="},
    @{old = "61÷2="; new = "92÷8="},
    @{old = "10÷8="; new = "57÷7="},
    @{old = "74÷7="; new = "94÷9="},
    @{old = "38÷4="; new = "58÷8="},
    @{old = "33÷7="; new = "57÷5="},
    @{old = "14÷8="; new = "69÷8="},
    @{old = "96÷3="; new = "54÷4="},
    @{old = "73÷8="; new = "96÷7="},
    @{old = "16÷3="; new = "72÷3="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
